$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C3").Value = 13194
$ws.Range("C4:C5").Value = 11670
$ws.Range("C6:C8").Value = 10799
$ws.Range("C9:C11").Value = 9551
$ws.Range("C12:C18").Value = 9459
$ws.Range("C19:C20").Value = 9451
$ws.Range("C21:C33").Value = 9020
$ws.Range("C34:C40").Value = 8341
$ws.Range("C41:C66").Value = 7769
$ws.Range("C67:C85").Value = 7318
$ws.Range("C86:C179").Value = 7293
